$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textForceCells = @("D4", "D5", "D6", "D9", "D11", "D12", "D13", "D14", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D38", "D39", "D40", "D42", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.422.04"
$ws.Range("E2").Value = "  -4.11%  "

$ws.Range("D3").Value = "3.277.26"
$ws.Range("E3").Value = "  -6.15%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "595.56"
$ws.Range("E5").Value = "  -3.44%  "

$ws.Range("D6").Value = "151.45"
$ws.Range("E6").Value = "  -10.73%  "

$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").Value = "3.268.50"
$ws.Range("E8").Value = "  -6.22%  "

$ws.Range("D9").Value = "0.545"
$ws.Range("E9").Value = "  -9.60%  "

$ws.Range("E10").Value = "  -12.24%  "

$ws.Range("D11").Value = "6.74"
$ws.Range("E11").Value = "  -6.37%  "

$ws.Range("D12").Value = "0.511"
$ws.Range("E12").Value = "  -11.18%  "

$ws.Range("D13").Value = "38.60"
$ws.Range("E13").Value = "  -14.72%  "

$ws.Range("D14").Value = "0.0000246"
$ws.Range("E14").Value = "  -9.42%  "

$ws.Range("D15").Value = "3.802.28"
$ws.Range("E15").Value = "  -6.46%  "

$ws.Range("D16").Value = "67.448.08"
$ws.Range("E16").Value = "  -4.28%  "

$ws.Range("D17").Value = "3.273.76"
$ws.Range("E17").Value = "  -6.69%  "

$ws.Range("D18").Value = "536.41"
$ws.Range("E18").Value = "  -10.17%  "

$ws.Range("E19").Value = "  -5.95%  "

$ws.Range("D20").Value = "7.21"
$ws.Range("E20").Value = "  -13.09%  "

$ws.Range("D21").Value = "15.18"
$ws.Range("E21").Value = "  -13.05%  "

$ws.Range("D22").Value = "0.764"
$ws.Range("E22").Value = "  -12.13%  "

$ws.Range("D23").Value = "7.88"
$ws.Range("E23").Value = "  -11.93%  "

$ws.Range("D24").Value = "86.09"
$ws.Range("E24").Value = "  -10.86%  "

$ws.Range("D25").Value = "13.58"
$ws.Range("E25").Value = "  -11.90%  "

$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("D27").Value = "3.27"
$ws.Range("E27").Value = "  -10.80%  "

$ws.Range("D28").Value = "2.17"
$ws.Range("E28").Value = "  -13.60%  "

$ws.Range("D29").Value = "29.49"
$ws.Range("E29").Value = "  -10.99%  "

$ws.Range("D30").Value = "8.07"
$ws.Range("E30").Value = "  -8.90%  "

$ws.Range("D31").Value = "2.69"
$ws.Range("E31").Value = "  -7.54%  "

$ws.Range("E32").Value = "  -8.66%  "

$ws.Range("D33").Value = "6.63"
$ws.Range("E33").Value = "  -17.22%  "

$ws.Range("D34").Value = "5.80"
$ws.Range("E34").Value = "  -12.99%  "

$ws.Range("D35").Value = "532.86"
$ws.Range("E35").Value = "  -9.74%  "

$ws.Range("E36").Value = "  -0.29%  "

$ws.Range("E37").Value = "  -7.43%  "

$ws.Range("D38").Value = "53.42"
$ws.Range("E38").Value = "  -5.67%  "

$ws.Range("D39").Value = "0.0862"
$ws.Range("E39").Value = "  -11.96%  "

$ws.Range("D40").Value = "9.06"
$ws.Range("E40").Value = "  -15.85%  "

$ws.Range("E41").Value = "  -10.98%  "

$ws.Range("D42").Value = "2.80"
$ws.Range("E42").Value = "  -15.53%  "

$ws.Range("D43").Value = "2.940.85"
$ws.Range("E43").Value = "  -11.03%  "

$ws.Range("D44").Value = "0.269"
$ws.Range("E44").Value = "  -11.64%  "

$ws.Range("D45").Value = "0.0₃0595"
$ws.Range("E45").Value = "  -16.28%  "

$ws.Range("E46").Value = "  -10.49%  "

$ws.Range("D47").Value = "26.92"
$ws.Range("E47").Value = "  -13.98%  "

$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "2.36"
$ws.Range("E48").Value = "  -16.03%  "

$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("D50").Value = "125.14"
$ws.Range("E50").Value = "  -6.10%  "

$ws.Range("D51").Value = "0.114"
$ws.Range("E51").Value = "  -11.25%  "
